# "zero out BLP for LCFS work"
#
# The BLP sheet's forecast row (row 2, columns B:AF) used to hold formulas
# that derived the BAU LCFS percentage trend line (interpolating between
# historic/forecast anchor points in B2 and L2, then a TREND() projection
# for M2:AF2). This change zeroes the whole row out, turning every one of
# those formula cells into a literal 0, while keeping the existing cell
# formatting intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BLP")

# Make BLP the active sheet/selection (matches the workbook view left behind
# by the author after performing this edit).
$ws.Activate() | Out-Null

# B2:L2 already carry style index 14 - just overwrite their formulas with 0.
$ws.Range("B2:L2").Value = 0

# M2:AF2 carried a different (percent-based) style that becomes unused once
# the TREND() formulas disappear. Zero the values, then copy L2's formatting
# onto them so the whole row ends up sharing the same plain style.
$ws.Range("M2:AF2").Value = 0
$ws.Range("L2").Copy() | Out-Null
$ws.Range("M2:AF2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Reflect the selection the author left behind after doing this edit.
$ws.Range("B2:AF2").Select() | Out-Null
